$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The whole table stores its values as text (see the workbook's
# numberStoredAsText ignoredError), so numeric/boolean-looking entries need
# a leading apostrophe to keep Excel from auto-converting them to a number
# or boolean. Plain words (Asset, Class 1, QQQ, the comparison operators)
# are unambiguous text already and don't need it.

# Update existing constraint row 8: Weight (G8) changes from 0 to 0.01
$ws.Range("G8").Value = "'0.01"

# Add a new constraint row 9
$ws.Range("A9").Value = "'8"
$ws.Range("B9").Value = "'FALSE"
$ws.Range("C9").Value = "Asset"
$ws.Range("D9").Value = "Class 1"
$ws.Range("E9").Value = "QQQ"
$ws.Range("F9").Value = ">"
$ws.Range("G9").Value = "'0.05"
